$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This handback-status report got a new localization round: the previously
# handed-back source file (UUID 879e77d3...) was regenerated under a new
# UUID (bd4c2053...) and a brand-new source file (UUID e8c2383c...) shows up
# as a content-duplicate of it. Update the three sheets (Overview, zh-cn,
# de-de) accordingly and grow each of their tables by one row.
# ---------------------------------------------------------------------------

$oldGuid = "879e77d3-ab9b-4ed8-8ae1-6a51bd3c5903"
$guid1   = "bd4c2053-d45e-480a-aca8-056d877ae46d"
$guid2   = "e8c2383c-8f44-4b42-8085-804ad6e80005"

$xlf1zh = "$guid1.5f7e76d22d1ae8563ba28a639bb9c5f37d651248.zh-cn.xlf"
$xlf1de = "$guid1.5f7e76d22d1ae8563ba28a639bb9c5f37d651248.de-de.xlf"
$xlf2zh = "$guid2.7c282366b5d539f6b88b14af72e916d04fdb65b8.zh-cn.xlf"
$xlf2de = "$guid2.7c282366b5d539f6b88b14af72e916d04fdb65b8.de-de.xlf"

$statusText = "Handed back: in sync with en-US"

# NB: this engine auto-coerces bare "True"/"False" into a real Boolean cell
# (t="b") and drops an assigned "" into no cell at all. The source data
# wants literal TEXT cells (t="s", including an explicit empty string), so
# force text interpretation the same way Excel's UI does for a leading
# apostrophe.
$trueText  = "'True"
$falseText = "'False"
$emptyText = "'"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Row 2 (existing row): file name / path get the regenerated GUID.
$ws.Range("A2").Value = "$guid1.md"
$ws.Range("C2").Value = ".md"
$ws.Range("E2").Value = $statusText
$ws.Range("F2").Value = $statusText
$ws.Range("G2").Value = "2016-08-29 23:04:41"

# Row 3 (new row): the content-duplicate file.
$ws.Range("A3").Value = "$guid2.md"
$ws.Range("C3").Value = ".md"
$ws.Range("E3").Value = $statusText
$ws.Range("F3").Value = $statusText
$ws.Range("G3").Value = "2016-08-29 23:04:41"

# Hyperlinks: this engine's Hyperlinks collection isn't per-cell addressable
# for updates, so clear them all and re-add in order (B2 keeps rId2, B3
# becomes the new rId3) matching the target ref ids.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c5eae0cdc6f89a1b511c27272241337a4fde9d6/e2e/$guid1.md", "", "", "e2e\$guid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c5eae0cdc6f89a1b511c27272241337a4fde9d6/e2e/$guid2.md", "", "", "e2e\$guid2.md") | Out-Null

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Row 2 (existing row) -> guid1
$ws.Range("A2").Value = "$guid1.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $statusText
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = $falseText
$ws.Range("G2").Value = $xlf1zh
$ws.Range("H2").Value = "2016-08-29 23:04:36"
$ws.Range("I2").Value = "$guid1.md"
$ws.Range("J2").Value = $xlf1zh
$ws.Range("K2").Value = "2016-08-29 23:04:53"
$ws.Range("L2").Value = $emptyText
$ws.Range("M2").Value = $trueText
$ws.Range("N2").Value = $emptyText
$ws.Range("O2").Value = $falseText
$ws.Range("P2").Value = $emptyText

# Row 3 (new row) -> guid2, content duplicate of row 2
$ws.Range("A3").Value = "$guid2.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $statusText
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = $trueText
$ws.Range("G3").Value = $xlf2zh
$ws.Range("H3").Value = "2016-08-29 23:04:36"
$ws.Range("I3").Value = "$guid2.md"
$ws.Range("J3").Value = $xlf2zh
$ws.Range("K3").Value = "2016-08-29 23:04:53"
$ws.Range("L3").Value = $emptyText
$ws.Range("M3").Value = $trueText
$ws.Range("N3").Value = $emptyText
$ws.Range("O3").Value = $falseText
$ws.Range("P3").Value = $emptyText

$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c5eae0cdc6f89a1b511c27272241337a4fde9d6/e2e/$guid1.md", "", "", "$guid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1ced93a3103dd645bef29fb41a04c395cfb5b973/e2e/$guid1.md", "", "", "$guid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c5eae0cdc6f89a1b511c27272241337a4fde9d6/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1ced93a3103dd645bef29fb41a04c395cfb5b973/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Row 2 (existing row) -> guid1
$ws.Range("A2").Value = "$guid1.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $statusText
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = $falseText
$ws.Range("G2").Value = $xlf1de
$ws.Range("H2").Value = "2016-08-29 23:04:41"
$ws.Range("I2").Value = "$guid1.md"
$ws.Range("J2").Value = $xlf1de
$ws.Range("K2").Value = "2016-08-29 23:05:00"
$ws.Range("L2").Value = $emptyText
$ws.Range("M2").Value = $trueText
$ws.Range("N2").Value = $emptyText
$ws.Range("O2").Value = $falseText
$ws.Range("P2").Value = $emptyText

# Row 3 (new row) -> guid2, content duplicate of row 2
$ws.Range("A3").Value = "$guid2.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $statusText
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = $trueText
$ws.Range("G3").Value = $xlf2de
$ws.Range("H3").Value = "2016-08-29 23:04:41"
$ws.Range("I3").Value = "$guid2.md"
$ws.Range("J3").Value = $xlf2de
$ws.Range("K3").Value = "2016-08-29 23:05:00"
$ws.Range("L3").Value = $emptyText
$ws.Range("M3").Value = $trueText
$ws.Range("N3").Value = $emptyText
$ws.Range("O3").Value = $falseText
$ws.Range("P3").Value = $emptyText

$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c5eae0cdc6f89a1b511c27272241337a4fde9d6/e2e/$guid1.md", "", "", "$guid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2680b7c099a06cde3e942c02d83ac3d9de77ca88/e2e/$guid1.md", "", "", "$guid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c5eae0cdc6f89a1b511c27272241337a4fde9d6/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2680b7c099a06cde3e942c02d83ac3d9de77ca88/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null
